$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for column G, rows 2-35, replacing the old Strike# values.
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 3
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 3
    34 = 1
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
